# Automatische test-sync: 2025-08-03 18:44:50
# Appends the new "Testmail #15" log entry to the Logs sheet, extends the
# conditional-formatting ranges to cover it, and bumps the matching
# "Planning / Afspraak" tally on the Dashboard sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 43

$logs.Cells.Item($newRow, 1).Value = "Leg dit even neer bij Koen."
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #15: Leg dit even neer bij Koen."
$logs.Cells.Item($newRow, 4).Value = "Planning / Afspraak"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-03 18:44:33"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Extend the conditional formatting that previously stopped at row 42 so it
# also covers the newly added row 43, for every column that has rules.
$dFc = $logs.Range("D2:D42").FormatConditions
for ($i = 1; $i -le $dFc.Count; $i++) {
    $dFc.Item($i).ModifyAppliesToRange($logs.Range("D2:D43"))
}

$gFc = $logs.Range("G2:G42").FormatConditions
for ($i = 1; $i -le $gFc.Count; $i++) {
    $gFc.Item($i).ModifyAppliesToRange($logs.Range("G2:G43"))
}

$hFc = $logs.Range("H2:H42").FormatConditions
for ($i = 1; $i -le $hFc.Count; $i++) {
    $hFc.Item($i).ModifyAppliesToRange($logs.Range("H2:H43"))
}

$iFc = $logs.Range("I2:I42").FormatConditions
for ($i = 1; $i -le $iFc.Count; $i++) {
    $iFc.Item($i).ModifyAppliesToRange($logs.Range("I2:I43"))
}

$jFc = $logs.Range("J2:J42").FormatConditions
for ($i = 1; $i -le $jFc.Count; $i++) {
    $jFc.Item($i).ModifyAppliesToRange($logs.Range("J2:J43"))
}

# Update the Dashboard summary count for "Planning / Afspraak" (row 3) to
# reflect the newly logged entry.
$dashboard.Range("B3").Value = 10
